# The worksheet carries sheet protection, so it must be unprotected before
# any cell content can be changed, then re-protected afterwards to restore
# the original protected state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect("D382")

# Update the "as of" date in the confidentiality disclosure (A10):
# 2021-03-26 -> 2021-03-29
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for each holding
$ws.Range("D2").Value = 0.4879947299476781
$ws.Range("E2").Value = -0.005110062893081802

$ws.Range("D3").Value = 0.3341208124159202
$ws.Range("E3").Value = -0.001964057743297554

$ws.Range("D4").Value = 0.09388668882713572
$ws.Range("E4").Value = -0.007715468573579409

$ws.Range("D5").Value = 0.05484589892504058
$ws.Range("E5").Value = 0.0005752416014725981

$ws.Range("D6").Value = 0.02915186988422543
$ws.Range("E6").Value = -0.01396103896103895

# Row 7 is the Total row; only the Percent Change figure changes
$ws.Range("E7").Value = -0.004249736875979626

# Restore sheet protection to match the original protected workbook
$ws.Protect("D382")
